# Fixed mappings for 07_08_2025.
# Update the account-name-to-number mappings on the "MS" sheet and
# move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MS")

# Map of cell -> new value (only the trailing account number changed)
$updates = @{
    "A3"  = "Cardor Alef SpA - 4165"
    "A4"  = "Cuenta Personal - Jaime - 4120"
    "A5"  = "NNW Capital SpA - 4105"
    "A6"  = "NNW II Capital SpA - 4103"
    "A7"  = "FNW Capital SpA - 4157"
    "A10" = "Alanseb LP - 4156"
    "A11" = "Cuenta Personal - Felipe - 4122"
    "A12" = "Cuenta Personal - Irene - 4162"
    "A13" = "Cuenta Personal - Jaime y Felipe - 4114"
    "A14" = "Cuenta Personal - Jaime y Natalia - 4118"
    "A15" = "Cuenta Personal - Jaime y Nicolas - 4121"
    "A16" = "Cuenta Personal - Natalia - 4108"
    "A17" = "Cuenta Personal - Nicolas - 4110"
    "A20" = "NNW Ventures LLC - 4107"
    "A21" = "Cuenta Personal - Jacques - 4155"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Activate the sheet and move selection to E10 (single cell, matching diff)
$ws.Activate()
$ws.Range("E10").Select()
